# Update external growth factor.
# Inserts 4 new DRI rows (TAZ 4860, 4861, 4902, 4903) into the sorted "DRIs"
# table, ahead of the existing TAZ 4957 row, shifting the remaining 5 rows
# (TAZ 4957, 5280, 5283, 5289, 5294) down from 95-99 to 99-103.
#
# NOTE: this is done as a value/format copy-down (not a structural
# Rows.Insert) so that the *other* worksheet (DRI_Increment), whose
# formulas point at DRIs!<row>, keeps referencing the same row numbers it
# always did -- exactly what the target file shows (formula text for
# DRI_Increment rows 95-99 is untouched; only the cached results change
# because the underlying DRIs numbers moved).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DRIs")
$ws2 = $wb.Worksheets.Item("DRI_Increment")

# ---------------------------------------------------------------------
# 1) Push the last 5 existing rows (95-99) down to (99-103), values +
#    formats together, so the table keeps growing from the bottom while
#    preserving ascending TAZ sort order once the new rows are written.
# ---------------------------------------------------------------------
$ws1.Range("A95:Q99").Copy()
$ws1.Range("A99:Q103").PasteSpecial(-4104)   # xlPasteAll
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Columns B:C of the 4 new rows (95-98) use a style that didn't exist
#    before: General number format with the workbook's red "variance"
#    font. Seed that by copying column A's plain/General format into B:C,
#    then recolor red so the engine reuses the existing red font object.
# ---------------------------------------------------------------------
$ws1.Range("A95:A98").Copy()
$ws1.Range("B95:C98").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("B95:C98").Font.Color = 255

# ---------------------------------------------------------------------
# 3) Write the new row values (columns A through Q).
# ---------------------------------------------------------------------
$newRows = @(
  @(4860, 671, 839, 1366, 1291, 2061, 1743, 2756, 2195, 3451, 2647, 4146, 3099, 4841, 3551, 5535, 4006),
  @(4861, 3293, 714, 3787, 2005, 4281, 3296, 4775, 4587, 5269, 5878, 5763, 7169, 6257, 8460, 6754, 9752),
  @(4902, 394, 1525, 667, 1632, 940, 1739, 1213, 1846, 1486, 1953, 1759, 2060, 2032, 2167, 2308, 2276),
  @(4903, 1305, 648, 1701, 791, 2097, 934, 2493, 1077, 2889, 1220, 3285, 1363, 3681, 1506, 4080, 1646)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = 95 + $i
  $rowVals = $newRows[$i]
  for ($c = 1; $c -le 17; $c++) {
    $ws1.Cells.Item($r, $c).Value2 = $rowVals[$c - 1]
  }
}

# ---------------------------------------------------------------------
# 4) The table's sort range / sort condition now covers the 4 extra rows.
# ---------------------------------------------------------------------
$ws1.Sort.SortFields.Clear()
$sf = $ws1.Sort.SortFields.Add($ws1.Range("A2:A103"))
$ws1.Sort.SetRange($ws1.Range("A2:Q103"))
$ws1.Sort.Header = -4142

# ---------------------------------------------------------------------
# 5) Leftover hidden "_xlnm._FilterDatabase" name for the DRIs sheet
#    (artifact of toggling AutoFilter on the data range at some point).
# ---------------------------------------------------------------------
$fdb = $ws1.Names.Add("_xlnm._FilterDatabase", "=DRIs!`$A`$1:`$Q`$99")
$fdb.Visible = $false

# ---------------------------------------------------------------------
# 6) View state: DRIs becomes the active tab/sheet, selection covers the
#    full refreshed table; DRI_Increment's selection/zoom settle elsewhere.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1:Q103").Select()
$ws1.Range("Q103").Activate()

$ws2.Select()
$ws2.Range("L12").Select()
$excel.ActiveWindow.Zoom = 100
